$d = $word.ActiveDocument
$sel = $word.Selection

# ---------------------------------------------------------------------------
# Edit 1: "Does selection affect homologs..." paragraph - replace the closing
# sentence of the explanatory text with a shorter replacement sentence.
# ---------------------------------------------------------------------------
$old1 = "First I will identify the homologous genes in the reference MAGs and SAGs using BLAST and clustering. Then I will mapping the metagenomic reads from the same lake back the reference genomes and determine the sequence-discrete populations each reference belongs to based on its coverage discontinuity. For each of these populations and their genes, I will calculate the relative abundance, call single nucleotide variants, and calculate the degree of selection."
$new1 = "First I will identify the homologous genes in the reference MAGs and SAGs using BLAST and clustering. Then I will mapping the metagenomic reads from the same lake back the reference genomes and determine the sequence-discrete populations each reference belongs to based on its coverage discontinuity. I will then calculate a metric of selection for within the population and between the populations."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Select()
    $sel.TypeText($new1)
}

# ---------------------------------------------------------------------------
# Edit 2: "Are there related sequence-discrete populations..." italic question
# - append an additional question to the end of the sentence.
# ---------------------------------------------------------------------------
$old2 = "Are there related sequence-discrete populations in TB and CB? How closely related are they? Do they share a common gene pool?"
$new2 = "Are there related sequence-discrete populations in TB and CB? How closely related are they? Do they share a common gene pool? Are there genes present in one lake but not the other for these closely related populations?"

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Select()
    $sel.TypeText($new2)
}

# ---------------------------------------------------------------------------
# Edit 3: "We can now bin MAG's from Crystal Bog..." paragraph - replace with
# the expanded/reworded paragraph, and drop the now-removed trailing space
# run that used to follow it.
# ---------------------------------------------------------------------------
$old3 = "We can now bin MAG's from Crystal Bog (CB), which is of similar location and trophic status to Trout Bog (TB). I will bin MAGs from the CB assemblies and use new techniques to get more MAGs from the TB assemblies. With these genomes, I will ask if there are related genomes and how similar they are across their genomes. For genomes that are very closely related (>95% nucleotide identity across their whole genomes), I will also investigate if the associated populations share a common gene pool between the two lakes. Using the cross mapping I will identify if there are regions or genes that are present in only one of the lakes and if the diversity of shared genes is different between the two lakes."
$new3 = "We now also have Crystal Bog metagenomes and can bin MAG's from Crystal Bog (CB), which is of similar location and trophic status to Trout Bog (TB), where our previous MAGs were from. I will bin MAGs from the CB assemblies and use new techniques to get more MAGs from the TB assemblies. With these genomes, I will search for very closely related genomes and quantify how similar they are across their genomes. For genomes that are very closely related (>95% nucleotide identity across their whole genomes), I will also investigate if the associated populations share a common gene pool between the two lakes. By mapping the metagenomes from one lake to MAGs from the other, I will identify if there are regions or genes that are present in only one of the lakes and if the diversity of shared genes is different between the two lakes. I will look for patterns among the shared or absent genes that may be explained by the different environments of the two lakes. I will also look to see if there is evidence of a barrier to recombination between these allopatric populations."

$rng3 = $d.Content
$found3 = $rng3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Select()
    $sel.TypeText($new3)

    # The paragraph used to end with the replaced text followed by a lone
    # " " run; that trailing space run was deleted in the target edit, so
    # remove the single space character that now immediately follows our
    # freshly typed paragraph text.
    $rng4 = $d.Content
    $found4 = $rng4.Find.Execute($new3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found4) {
        $trailing = $d.Range($rng4.End, $rng4.End + 1)
        if ($trailing.Text -eq " ") {
            $trailing.Delete()
        }
    }
}
